# Update the "Clasificación" (standings) table with the latest match results.
# Each player's PJ (games played) goes up by one, and several derived
# columns (TD/TP/DT/V1/V2/VC/PTS) are refreshed to reflect the new totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - David
$ws.Range("C2").Value = 8
$ws.Range("F2").Value = 3
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 2

# Row 3 - Pedro
$ws.Range("C3").Value = 8
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 9
$ws.Range("I3").Value = -4

# Row 4 - Adonay
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 8
$ws.Range("G4").Value = 11
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 10
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 27

# Row 5 - Richard
$ws.Range("C5").Value = 8
$ws.Range("F5").Value = 7
$ws.Range("H5").Value = 8
$ws.Range("I5").Value = -8

# Row 6 - Iván
$ws.Range("C6").Value = 8
$ws.Range("D6").Value = 5
$ws.Range("G6").Value = 9
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 3
$ws.Range("M6").Value = 19

# Row 7 - Nico
$ws.Range("C7").Value = 8
$ws.Range("D7").Value = 3
$ws.Range("G7").Value = 3
$ws.Range("I7").Value = -1
$ws.Range("J7").Value = 3
$ws.Range("M7").Value = 10

# Row 8 - Nicolás
$ws.Range("C8").Value = 8
$ws.Range("D8").Value = 3
$ws.Range("G8").Value = 3
$ws.Range("I8").Value = -4
$ws.Range("J8").Value = 3
$ws.Range("M8").Value = 9

# Row 9 - Vicente
$ws.Range("C9").Value = 8
$ws.Range("F9").Value = 4
$ws.Range("H9").Value = 4
$ws.Range("I9").Value = 0
